$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "asdf"
$ws.Range("B2").Value = "asdfasdf"
$ws.Range("C2").Value = 9112868830
$ws.Range("D2").Value = "erfangg"
$ws.Range("E2").Value = "Erfan@gmail.com"
$ws.Range("F2").Value = "Erfan@123"
$ws.Range("G2").Value = "Tehran"
$ws.Range("H2").Value = "2005/February/26"

$ws.Range("A3").Value = "Erfan"
$ws.Range("B3").Value = "Ghasemian"
$ws.Range("C3").Value = "'09112868820"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "erfangh"
$ws.Range("E3").Value = "erfan.ghasemian40@gmail.com"
$ws.Range("F3").Value = "Erfan@123"
$ws.Range("G3").Value = "Tehran"
$ws.Range("H3").Value = "2005/January/1"
